$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at position 667 (pushes existing 667:750 down to 670:753)
$ws.Rows("667:669").Insert()

# Row 667 - new weekly data point (Tomate, Larga vida, Primera)
$ws.Range("A667").Value = 9
$ws.Range("B667").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C667").Value = "Metropolitana"
$ws.Range("D667").Value = 44474
$ws.Range("E667").Value = 13
$ws.Range("F667").Value = 100112020
$ws.Range("G667").Value = "Tomate"
$ws.Range("H667").Value = "Larga vida"
$ws.Range("I667").Value = "Primera"
$ws.Range("J667").Value = 250
$ws.Range("K667").Value = 23000
$ws.Range("L667").Value = 24000
$ws.Range("M667").Value = 23500
$ws.Range("N667").Value = "`$/bandeja 20 kilos"
$ws.Range("O667").Value = "Región de Arica y Parinacota"
$ws.Range("P667").Value = 1175
$ws.Range("Q667").Value = 20
$ws.Range("R667").Value = "Hortaliza"

# Row 668 - new weekly data point (Tomate, Larga vida, Segunda)
$ws.Range("A668").Value = 9
$ws.Range("B668").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C668").Value = "Metropolitana"
$ws.Range("D668").Value = 44474
$ws.Range("E668").Value = 13
$ws.Range("F668").Value = 100112020
$ws.Range("G668").Value = "Tomate"
$ws.Range("H668").Value = "Larga vida"
$ws.Range("I668").Value = "Segunda"
$ws.Range("J668").Value = 160
$ws.Range("K668").Value = 21000
$ws.Range("L668").Value = 22000
$ws.Range("M668").Value = 21500
$ws.Range("N668").Value = "`$/bandeja 20 kilos"
$ws.Range("O668").Value = "Región de Arica y Parinacota"
$ws.Range("P668").Value = 1075
$ws.Range("Q668").Value = 20
$ws.Range("R668").Value = "Hortaliza"

# Row 669 - new weekly data point (Tomate, Larga vida, Tercera)
$ws.Range("A669").Value = 9
$ws.Range("B669").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C669").Value = "Metropolitana"
$ws.Range("D669").Value = 44474
$ws.Range("E669").Value = 13
$ws.Range("F669").Value = 100112020
$ws.Range("G669").Value = "Tomate"
$ws.Range("H669").Value = "Larga vida"
$ws.Range("I669").Value = "Tercera"
$ws.Range("J669").Value = 97
$ws.Range("K669").Value = 19000
$ws.Range("L669").Value = 20000
$ws.Range("M669").Value = 19495
$ws.Range("N669").Value = "`$/bandeja 20 kilos"
$ws.Range("O669").Value = "Región de Arica y Parinacota"
$ws.Range("P669").Value = 975
$ws.Range("Q669").Value = 20
$ws.Range("R669").Value = "Hortaliza"
